$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.02325754048051
$ws.Range("D2").Value = 1.033211150923203
$ws.Range("E2").Value = 0.9926147277508489
$ws.Range("F2").Value = 1.041933255036759
$ws.Range("I2").Value = 1.033649761585882
$ws.Range("J2").Value = 1.028438773386392
$ws.Range("K2").Value = 1.036014153706946
$ws.Range("L2").Value = 0.9955398523336033
$ws.Range("M2").Value = 1.044711369327488
$ws.Range("N2").Value = 1.013554452336197
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.024056269950871
$ws.Range("D3").Value = 1.033813372346365
$ws.Range("E3").Value = 0.9936372048519304
$ws.Range("F3").Value = 1.042673478669354
$ws.Range("I3").Value = 1.033780111734292
$ws.Range("J3").Value = 1.028876685618918
$ws.Range("K3").Value = 1.036425860494977
$ws.Range("L3").Value = 0.9963617723202692
$ws.Range("M3").Value = 1.045262506603683
$ws.Range("N3").Value = 1.013698688352498
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.024573448591761
$ws.Range("D4").Value = 1.034203004237257
$ws.Range("E4").Value = 0.9942998659930995
$ws.Range("F4").Value = 1.043152736510115
$ws.Range("I4").Value = 1.033862755391051
$ws.Range("J4").Value = 1.029159733532211
$ws.Range("K4").Value = 1.036691521863635
$ws.Range("L4").Value = 0.9968940712668345
$ws.Range("M4").Value = 1.04561874623777
$ws.Range("N4").Value = 1.013791909532652
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.02479095190174
$ws.Range("D5").Value = 1.034366792600601
$ws.Range("E5").Value = 0.9945786998346017
$ws.Range("F5").Value = 1.04335428228484
$ws.Range("I5").Value = 1.033897090459359
$ws.Range("J5").Value = 1.029278651202363
$ws.Range("K5").Value = 1.036803027149895
$ws.Range("L5").Value = 0.997117960005301
$ws.Range("M5").Value = 1.045768415901112
$ws.Range("N5").Value = 1.013831073134036
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.024827476388193
$ws.Range("D6").Value = 1.034394292574931
$ws.Range("E6").Value = 0.9946255319796338
$ws.Range("F6").Value = 1.0433881265025
$ws.Range("I6").Value = 1.033902831496118
$ws.Range("J6").Value = 1.029298613526551
$ws.Range("K6").Value = 1.036821738852573
$ws.Range("L6").Value = 0.9971555583673453
$ws.Range("M6").Value = 1.045793540562661
$ws.Range("N6").Value = 1.013837647303157
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.024576354563921
$ws.Range("D7").Value = 1.034205192838484
$ws.Range("E7").Value = 0.9943035907982488
$ws.Range("F7").Value = 1.043155429318846
$ws.Range("I7").Value = 1.033863215783561
$ws.Range("J7").Value = 1.029161322814945
$ws.Range("K7").Value = 1.036693012506985
$ws.Range("L7").Value = 0.9968970624462087
$ws.Range("M7").Value = 1.045620746500693
$ws.Range("N7").Value = 1.013792432943628
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.023527401904019
$ws.Range("D8").Value = 1.033414682748088
$ws.Range("E8").Value = 0.9929600610674301
$ws.Range("F8").Value = 1.042183357068992
$ws.Range("I8").Value = 1.033694165783241
$ws.Range("J8").Value = 1.028586831392673
$ws.Range("K8").Value = 1.036153444366914
$ws.Range("L8").Value = 0.995817528259106
$ws.Range("M8").Value = 1.044897706992932
$ws.Range("N8").Value = 1.013603219908144
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.021681749069477
$ws.Range("D9").Value = 1.032021439603961
$ws.Range("E9").Value = 0.9906006454969559
$ws.Range("F9").Value = 1.040472707775807
$ws.Range("I9").Value = 1.033383287891599
$ws.Range("J9").Value = 1.027572183459717
$ws.Range("K9").Value = 1.035197050814366
$ws.Range("L9").Value = 0.9939188001724441
$ws.Range("M9").Value = 1.043620763707908
$ws.Range("N9").Value = 1.013268985997175
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.020453255374702
$ws.Range("D10").Value = 1.031092549711188
$ws.Range("E10").Value = 0.989033133672735
$ws.Range("F10").Value = 1.03933392471328
$ws.Range("I10").Value = 1.033167358497147
$ws.Range("J10").Value = 1.026894269204073
$ws.Range("K10").Value = 1.034555774147882
$ws.Range("L10").Value = 0.9926553831429383
$ws.Range("M10").Value = 1.04276765370544
$ws.Range("N10").Value = 1.013045641007288
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.019921785933563
$ws.Range("D11").Value = 1.030690340025126
$ws.Range("E11").Value = 0.988355674866747
$ws.Range("F11").Value = 1.038841236913585
$ws.Range("I11").Value = 1.033071812289958
$ws.Range("J11").Value = 1.026600389620633
$ws.Range("K11").Value = 1.034277239631596
$ws.Range("L11").Value = 0.9921088820399291
$ws.Range("M11").Value = 1.042397839390411
$ws.Range("N11").Value = 1.012948811819041
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.019724447805681
$ws.Range("D12").Value = 1.030540944268971
$ws.Range("E12").Value = 0.9881042295826724
$ws.Range("F12").Value = 1.038658294879768
$ws.Range("I12").Value = 1.033036015503644
$ws.Range("J12").Value = 1.026491180148354
$ws.Range("K12").Value = 1.034173652152364
$ws.Range("L12").Value = 0.9919059725120875
$ws.Range("M12").Value = 1.042260413751365
$ws.Range("N12").Value = 1.01291282766309
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.019766774160247
$ws.Range("D13").Value = 1.030572990006325
$ws.Range("E13").Value = 0.9881581567098651
$ws.Range("F13").Value = 1.038697533629808
$ws.Range("I13").Value = 1.033043707898576
$ws.Range("J13").Value = 1.026514608170383
$ws.Range("K13").Value = 1.034195877757987
$ws.Range("L13").Value = 0.9919494934313052
$ws.Range("M13").Value = 1.042289894710398
$ws.Range("N13").Value = 1.012920547170767
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.019905472397469
$ws.Range("D14").Value = 1.030677990856272
$ws.Range("E14").Value = 0.9883348863814464
$ws.Range("F14").Value = 1.038826113546995
$ws.Range("I14").Value = 1.033068859567112
$ws.Range("J14").Value = 1.026591363334036
$ws.Range("K14").Value = 1.034268679647838
$ws.Range("L14").Value = 0.9920921077337197
$ws.Range("M14").Value = 1.042386480960771
$ws.Range("N14").Value = 1.012945837711591
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.019990938721208
$ws.Range("D15").Value = 1.030742685778917
$ws.Range("E15").Value = 0.9884438009545853
$ws.Range("F15").Value = 1.03890534444216
$ws.Range("I15").Value = 1.033084315730838
$ws.Range("J15").Value = 1.026638648202046
$ws.Range("K15").Value = 1.034313518462607
$ws.Range("L15").Value = 0.9921799884222134
$ws.Range("M15").Value = 1.042445983021534
$ws.Range("N15").Value = 1.012961417744937
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.020488537434878
$ws.Range("D16").Value = 1.031119243332685
$ws.Range("E16").Value = 0.9890781214508737
$ws.Range("F16").Value = 1.039366631659267
$ws.Range("I16").Value = 1.033173656527554
$ws.Range("J16").Value = 1.026913766009495
$ws.Range("K16").Value = 1.034574241625708
$ws.Range("L16").Value = 0.9926916645766087
$ws.Range("M16").Value = 1.042792188533312
$ws.Range("N16").Value = 1.01305206476536
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.020800796722745
$ws.Range("D17").Value = 1.0313554508703
$ws.Range("E17").Value = 0.989476357848556
$ws.Range("F17").Value = 1.039656096967283
$ws.Range("I17").Value = 1.033229150165595
$ws.Range("J17").Value = 1.027086250509484
$ws.Range("K17").Value = 1.034737557843372
$ws.Range("L17").Value = 0.9930127773699352
$ws.Range("M17").Value = 1.04300924503735
$ws.Range("N17").Value = 1.013108893611624
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.020982978137111
$ws.Range("D18").Value = 1.031493227155956
$ws.Range("E18").Value = 0.9897087662937556
$ws.Range("F18").Value = 1.039824976868666
$ws.Range("I18").Value = 1.033261321064228
$ws.Range("J18").Value = 1.027186825077786
$ws.Range("K18").Value = 1.034832734596681
$ws.Range("L18").Value = 0.9932001317071769
$ws.Range("M18").Value = 1.043135810548276
$ws.Range("N18").Value = 1.013142029372488
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.021045105017211
$ws.Range("D19").Value = 1.031540205321186
$ws.Range("E19").Value = 0.9897880325774034
$ws.Range("F19").Value = 1.039882567216875
$ws.Range("I19").Value = 1.033272256964554
$ws.Range("J19").Value = 1.027221112816411
$ws.Range("K19").Value = 1.034865173315189
$ws.Range("L19").Value = 0.9932640239640975
$ws.Range("M19").Value = 1.04317895929921
$ws.Range("N19").Value = 1.013153325839082
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.020767289501035
$ws.Range("D20").Value = 1.031330107978782
$ws.Range("E20").Value = 0.9894336180360679
$ws.Range("F20").Value = 1.039625035959005
$ws.Range("I20").Value = 1.03322321665561
$ws.Range("J20").Value = 1.027067747930193
$ws.Range("K20").Value = 1.034720044116041
$ws.Range("L20").Value = 0.9929783193494215
$ws.Range("M20").Value = 1.04298596104039
$ws.Range("N20").Value = 1.013102797605048
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.019864627186811
$ws.Range("D21").Value = 1.030647070629205
$ws.Range("E21").Value = 0.9882828385668249
$ws.Range("F21").Value = 1.038788248178131
$ws.Range("I21").Value = 1.033061461485855
$ws.Range("J21").Value = 1.026568762207662
$ws.Range("K21").Value = 1.034247244809752
$ws.Range("L21").Value = 0.9920501090198102
$ws.Range("M21").Value = 1.042358040363266
$ws.Range("N21").Value = 1.012938390756754
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.019297512418504
$ws.Range("D22").Value = 1.03021763527484
$ws.Range("E22").Value = 0.9875604150241495
$ws.Range("F22").Value = 1.038262498703266
$ws.Range("I22").Value = 1.032957985741129
$ws.Range("J22").Value = 1.026254744273091
$ws.Range("K22").Value = 1.033949240988709
$ws.Range("L22").Value = 0.9914670000341481
$ws.Range("M22").Value = 1.041962894651605
$ws.Range("N22").Value = 1.012834920707088
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.019598109846459
$ws.Range("D23").Value = 1.03044528478475
$ws.Range("E23").Value = 0.9879432794643023
$ws.Range("F23").Value = 1.03854117245447
$ws.Range("I23").Value = 1.033013008025165
$ws.Range("J23").Value = 1.026421237777704
$ws.Range("K23").Value = 1.034107287740324
$ws.Range("L23").Value = 0.991776070289318
$ws.Range("M23").Value = 1.042172401177815
$ws.Range("N23").Value = 1.012889781557988
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.020782429830735
$ws.Range("D24").Value = 1.031341559339512
$ws.Range("E24").Value = 0.9894529299347244
$ws.Range("F24").Value = 1.039639070969281
$ws.Range("I24").Value = 1.033225898364023
$ws.Range("J24").Value = 1.0270761085506
$ws.Range("K24").Value = 1.034727958071665
$ws.Range("L24").Value = 0.9929938892766442
$ws.Range("M24").Value = 1.042996482199787
$ws.Range("N24").Value = 1.01310555216374
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.022158559829213
$ws.Range("D25").Value = 1.032381645583551
$ws.Range("E25").Value = 0.9912096547607049
$ws.Range("F25").Value = 1.040914670038429
$ws.Range("I25").Value = 1.033465190854522
$ws.Range("J25").Value = 1.027834761228282
$ws.Range("K25").Value = 1.035444956678367
$ws.Range("L25").Value = 0.9944092447426414
$ws.Range("M25").Value = 1.043951211649241
$ws.Range("N25").Value = 1.013355487481978
